# Fix the "Sections" seed worksheet: several rows in column C hold the
# string "vid " (with a trailing space) instead of "vid" like the other
# rows. Normalize them so all rows use the same "vid" string, then leave
# the final selection on C20 (the last cell touched).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 6; $r -le 20; $r++) {
    $ws.Cells.Item($r, 3).Value = "vid"
}

$ws.Range("C20").Select()
